$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 785, shifting existing rows (785-826) down to (786-827)
$ws.Rows.Item(785).Insert()

# Populate the new row 785 with the inserted data
# Force column A to be stored as plain text (matching the other date-as-text cells
# in this sheet), not auto-converted to an Excel date serial number. Temporarily
# apply a Text number format so the value is kept literal, then restore the
# "Normal" cell style so no stray number-format override is left behind.
$ws.Cells.Item(785, 1).NumberFormat = "@"
$ws.Cells.Item(785, 1).Value = "2026/02/05"
$ws.Cells.Item(785, 1).Style = "Normal"
$ws.Cells.Item(785, 2).Value = "木"
$ws.Cells.Item(785, 3).Value = 14
$ws.Cells.Item(785, 4).Value = 201
